# Scheduled-runner refresh: pull the latest Market Board averages into each
# crafting-job sheet and recompute the Leve profit columns (H:N).
# currentAveragePrice / currentAveragePriceNQ / currentAveragePriceHQ (H,I,J)
# and the derived LevePriceNQ/HQ + LeveProfitNQ/HQ (K,L,M,N) move together;
# a handful of rows gain/lose the HQ profit cell entirely when HQ pricing
# becomes (un)available for that item.

$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 522.8570999999999
$ws.Range("I4").Value = 522.8570999999999
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 522.8570999999999
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -408.8570999999999
$ws.Range("N4").ClearContents()
$ws.Range("H8").Value = 408.875
$ws.Range("I8").Value = 457.57144
$ws.Range("K8").Value = 1372.71432
$ws.Range("M8").Value = -1233.71432
$ws.Range("H15").Value = 1756.1111
$ws.Range("I15").Value = 1756.1111
$ws.Range("K15").Value = 5268.3333
$ws.Range("M15").Value = -5099.3333
$ws.Range("H17").Value = 713.1795
$ws.Range("I17").Value = 972.5
$ws.Range("J17").Value = 683.54285
$ws.Range("K17").Value = 2917.5
$ws.Range("L17").Value = 2050.62855
$ws.Range("M17").Value = -2749.5
$ws.Range("N17").Value = -2386.62855
$ws.Range("H55").Value = 328.5
$ws.Range("I55").Value = 237.88889
$ws.Range("K55").Value = 237.88889
$ws.Range("M55").Value = -23.88889
$ws.Range("H70").Value = 2328.5
$ws.Range("I70").Value = 1150
$ws.Range("J70").Value = 3114.1667
$ws.Range("K70").Value = 3450
$ws.Range("L70").Value = 9342.500100000001
$ws.Range("M70").Value = -3180
$ws.Range("N70").Value = -9882.500100000001
$ws.Range("H73").Value = 2328.5
$ws.Range("I73").Value = 1150
$ws.Range("J73").Value = 3114.1667
$ws.Range("K73").Value = 3450
$ws.Range("L73").Value = 9342.500100000001
$ws.Range("M73").Value = -2514
$ws.Range("N73").Value = -11214.5001
$ws.Range("H76").Value = 3000
$ws.Range("J76").Value = 3000
$ws.Range("L76").Value = 3000
$ws.Range("N76").Value = -3630
$ws.Range("H79").Value = 3000
$ws.Range("J79").Value = 3000
$ws.Range("L79").Value = 3000
$ws.Range("N79").Value = -5184
$ws.Range("H98").Value = 1431
$ws.Range("I98").Value = 1416.8334
$ws.Range("J98").Value = 1686
$ws.Range("K98").Value = 1416.8334
$ws.Range("L98").Value = 1686
$ws.Range("M98").Value = 81.16660000000002
$ws.Range("N98").Value = -4682
$ws.Range("H115").Value = 799.6667
$ws.Range("I115").Value = 500
$ws.Range("K115").Value = 1500
$ws.Range("M115").Value = 67
$ws.Range("H122").Value = 1431
$ws.Range("I122").Value = 1416.8334
$ws.Range("J122").Value = 1686
$ws.Range("K122").Value = 4250.5002
$ws.Range("L122").Value = 5058
$ws.Range("M122").Value = -1800.5002
$ws.Range("N122").Value = -9958
$ws.Range("H125").Value = 7145006.5
$ws.Range("I125").Value = 2818.4285
$ws.Range("J125").Value = 14287195
$ws.Range("K125").Value = 25365.8565
$ws.Range("L125").Value = 128584755
$ws.Range("M125").Value = -22905.8565
$ws.Range("N125").Value = -128589675
$ws.Range("H137").Value = 17243824
$ws.Range("I137").Value = 33335274
$ws.Range("J137").Value = 2985.4285
$ws.Range("K137").Value = 100005822
$ws.Range("L137").Value = 8956.2855
$ws.Range("M137").Value = -100003272
$ws.Range("N137").Value = -14056.2855

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 17604.54
$ws.Range("I32").Value = 18445.322
$ws.Range("K32").Value = 18445.322
$ws.Range("M32").Value = -18158.322
$ws.Range("H33").Value = 35000
$ws.Range("I33").Value = 0
$ws.Range("J33").Value = 35000
$ws.Range("K33").Value = 0
$ws.Range("L33").Value = 35000
$ws.Range("M33").ClearContents()
$ws.Range("N33").Value = -35658
$ws.Range("H124").Value = 47998.668
$ws.Range("J124").Value = 47998.668
$ws.Range("L124").Value = 47998.668
$ws.Range("N124").Value = -57818.668

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1704.4286
$ws.Range("I20").Value = 1589
$ws.Range("K20").Value = 1589
$ws.Range("M20").Value = -1342
$ws.Range("H86").Value = 1563.25
$ws.Range("I86").Value = 1615.3334
$ws.Range("J86").Value = 1407
$ws.Range("K86").Value = 1615.3334
$ws.Range("L86").Value = 1407
$ws.Range("M86").Value = -492.3334
$ws.Range("N86").Value = -3653
$ws.Range("H89").Value = 1563.25
$ws.Range("I89").Value = 1615.3334
$ws.Range("J89").Value = 1407
$ws.Range("K89").Value = 8076.666999999999
$ws.Range("L89").Value = 7035
$ws.Range("M89").Value = -2460.666999999999
$ws.Range("N89").Value = -18267
$ws.Range("H107").Value = 2420.4736
$ws.Range("I107").Value = 1667.8214
$ws.Range("J107").Value = 4527.9
$ws.Range("K107").Value = 1667.8214
$ws.Range("L107").Value = 4527.9
$ws.Range("M107").Value = 252.1786
$ws.Range("N107").Value = -8367.9

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 3281.1428
$ws.Range("I62").Value = 2756
$ws.Range("J62").Value = 3675
$ws.Range("K62").Value = 2756
$ws.Range("L62").Value = 3675
$ws.Range("M62").Value = -2132
$ws.Range("N62").Value = -4923
$ws.Range("H65").Value = 3281.1428
$ws.Range("I65").Value = 2756
$ws.Range("J65").Value = 3675
$ws.Range("K65").Value = 13780
$ws.Range("L65").Value = 18375
$ws.Range("M65").Value = -10660
$ws.Range("N65").Value = -24615
$ws.Range("H94").Value = 1867
$ws.Range("I94").Value = 833.6667
$ws.Range("J94").Value = 2088.4285
$ws.Range("K94").Value = 833.6667
$ws.Range("L94").Value = 2088.4285
$ws.Range("M94").Value = -382.6667
$ws.Range("N94").Value = -2990.4285
$ws.Range("H132").Value = 40202750
$ws.Range("I132").Value = 47620896
$ws.Range("K132").Value = 142862688
$ws.Range("M132").Value = -142860158

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 77252170
$ws.Range("I4").Value = 132286290
$ws.Range("J4").Value = 204399.8
$ws.Range("K4").Value = 396858870
$ws.Range("L4").Value = 613199.3999999999
$ws.Range("M4").Value = -396858758
$ws.Range("N4").Value = -613423.3999999999
$ws.Range("H132").Value = 592.1429000000001
$ws.Range("I132").Value = 579.1667
$ws.Range("K132").Value = 5212.5003
$ws.Range("M132").Value = -2682.5003

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4589.4
$ws.Range("I70").Value = 4299.5
$ws.Range("J70").Value = 4782.6665
$ws.Range("K70").Value = 4299.5
$ws.Range("L70").Value = 4782.6665
$ws.Range("M70").Value = -4029.5
$ws.Range("N70").Value = -5322.6665
$ws.Range("H73").Value = 4589.4
$ws.Range("I73").Value = 4299.5
$ws.Range("J73").Value = 4782.6665
$ws.Range("K73").Value = 4299.5
$ws.Range("L73").Value = 4782.6665
$ws.Range("M73").Value = -3363.5
$ws.Range("N73").Value = -6654.6665
$ws.Range("H80").Value = 345273.8
$ws.Range("I80").Value = 572567.4399999999
$ws.Range("K80").Value = 572567.4399999999
$ws.Range("M80").Value = -571569.4399999999
$ws.Range("H83").Value = 345273.8
$ws.Range("I83").Value = 572567.4399999999
$ws.Range("K83").Value = 2862837.2
$ws.Range("M83").Value = -2857845.2
$ws.Range("H107").Value = 173672.5
$ws.Range("J107").Value = 0
$ws.Range("L107").Value = 0
$ws.Range("N107").ClearContents()

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H36").Value = 0
$ws.Range("J36").Value = 0
$ws.Range("L36").Value = 0
$ws.Range("N36").ClearContents()
$ws.Range("H40").Value = 3561
$ws.Range("I40").Value = 3396.375
$ws.Range("K40").Value = 3396.375
$ws.Range("M40").Value = -3260.375
$ws.Range("H46").Value = 6028.857
$ws.Range("I46").Value = 8489.666999999999
$ws.Range("J46").Value = 1599.4
$ws.Range("K46").Value = 8489.666999999999
$ws.Range("L46").Value = 1599.4
$ws.Range("M46").Value = -8301.666999999999
$ws.Range("N46").Value = -1975.4
$ws.Range("H124").Value = 79500
$ws.Range("J124").Value = 79500
$ws.Range("L124").Value = 79500
$ws.Range("N124").Value = -89320

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 125320150
$ws.Range("I4").Value = 1259997.5
$ws.Range("J4").Value = 166673540
$ws.Range("K4").Value = 1259997.5
$ws.Range("L4").Value = 166673540
$ws.Range("M4").Value = -1259884.5
$ws.Range("N4").Value = -166673766
$ws.Range("H123").Value = 74994.664
$ws.Range("J123").Value = 74994.664
$ws.Range("L123").Value = 74994.664
$ws.Range("N123").Value = -84794.664
